$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column I: CSS rule generator, built from H (nth-child index), G (X offset), F (Y offset)
# Row 3 first (creates the "master" formula), then I4:I14 as one range assignment so the
# engine groups it into a shared formula like the target workbook.
$ws.Range("I3").Formula = '=".filter-btn.open .panelTool:nth-child("&H3&") {   transform: translate("&INT(G3)&"px, "&INT(F3)&"px); }"'
$ws.Range("I4:I14").Formula = '=".filter-btn.open .panelTool:nth-child("&H4&") {   transform: translate("&INT(G4)&"px, "&INT(F4)&"px); }"'

# Column J: only J3 gets an (empty-string) helper formula
$ws.Range("J3").Formula = '=""'

# Widen column I so the generated CSS text is fully visible
$ws.Columns.Item(9).ColumnWidth = 87.5

# Select the newly generated CSS column, mirroring the author's final selection
$null = $ws.Range("I3:I13").Select()
